# This script updates the "Price" (D) and "Volume(1h)" (E) columns of the
# cryptos worksheet with refreshed values from the GitHub Actions data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dUpdates = [ordered]@{
    "D2" = "30.065.80"
    "D3" = "1.892.06"
    "D4" = "1.000"
    "D5" = "249.08"
    "D6" = "1.000"
    "D7" = "0.5005"
    "D8" = "45.82"
    "D9" = "0.2868"
    "D10" = "0.06565"
    "D11" = "1.883.13"
    "D12" = "17.24"
    "D13" = "0.07230"
    "D14" = "0.6690"
    "D15" = "85.14"
    "D16" = "4.830"
    "D17" = "30.077.50"
    "D18" = "0.9999"
    "D19" = "12.91"
    "D20" = "0.000007546"
    "D21" = "1.000"
    "D22" = "2.128.67"
    "D23" = "4.778"
    "D24" = "5.555"
    "D25" = "9.053"
    "D26" = "145.16"
    "D27" = "135.90"
    "D28" = "16.80"
    "D29" = "1.957"
    "D30" = "1.371"
    "D31" = "4.202"
    "D32" = "0.08679"
    "D33" = "3.932"
    "D34" = "0.05047"
    "D36" = "0.6938"
    "D37" = "2.686"
    "D38" = "2.282"
    "D39" = "2.773"
    "D41" = "0.01644"
    "D42" = "6.048"
    "D43" = "105.90"
    "D44" = "1.000"
    "D45" = "0.4216"
    "D46" = "7.453"
    "D48" = "0.05666"
    "D49" = "32.59"
    "D50" = "8.284"
    "D51" = "0.3724"
}

$eUpdates = [ordered]@{
    "E2" = "  +6.14%  "
    "E3" = "  +5.69%  "
    "E4" = "  +0.13%  "
    "E5" = "  +0.63%  "
    "E6" = "  +0.14%  "
    "E7" = "  +1.65%  "
    "E8" = "  +8.19%  "
    "E9" = "  +6.27%  "
    "E10" = "  +4.26%  "
    "E11" = "  +5.12%  "
    "E12" = "  +3.69%  "
    "E13" = "  +2.54%  "
    "E14" = "  +6.10%  "
    "E15" = "  +6.13%  "
    "E16" = "  +3.57%  "
    "E17" = "  +6.26%  "
    "E18" = "  +0.08%  "
    "E19" = "  +6.77%  "
    "E20" = "  +3.88%  "
    "E21" = "  +0.16%  "
    "E22" = "  +5.58%  "
    "E23" = "  +4.58%  "
    "E24" = "  +5.63%  "
    "E25" = "  +3.20%  "
    "E26" = "  +2.31%  "
    "E27" = "  +23.61%  "
    "E28" = "  +6.34%  "
    "E29" = "  +5.22%  "
    "E30" = "  -1.37%  "
    "E31" = "  +0.60%  "
    "E32" = "  +4.79%  "
    "E33" = "  +3.99%  "
    "E34" = "  +3.13%  "
    "E35" = "  +5.20%  "
    "E36" = "  +5.03%  "
    "E37" = "  +2.79%  "
    "E38" = "  +9.69%  "
    "E39" = "  +5.92%  "
    "E40" = "  +1.78%  "
    "E41" = "  +5.49%  "
    "E42" = "  +2.04%  "
    "E43" = "  +6.16%  "
    "E44" = "  +0.17%  "
    "E45" = "  +5.22%  "
    "E46" = "  +3.37%  "
    "E47" = "  +3.65%  "
    "E49" = "  +5.75%  "
    "E50" = "  +2.59%  "
    "E51" = "  +6.31%  "
}

# Keep a reference style (from an untouched data cell) so that cells we
# rewrite as text keep the workbook's original (default) cell style
# instead of picking up a new "quote-prefixed"/text style.
$refStyle = $ws.Range("B2").Style

foreach ($addr in $dUpdates.Keys) {
    $cell = $ws.Range($addr)
    # Force the cell to Text format before assigning so values that look
    # numeric (e.g. "1.000", "105.90") are preserved exactly as typed
    # instead of being normalized into a number.
    $cell.NumberFormat = "@"
    $cell.Value = $dUpdates[$addr]
    $cell.Style = $refStyle
}

foreach ($addr in $eUpdates.Keys) {
    $ws.Range($addr).Value = $eUpdates[$addr]
}

Write-Output "Updated $($dUpdates.Count) Price cells and $($eUpdates.Count) Volume(1h) cells"
